$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 62: B62 changes from 2 to 2.5
$ws.Range("B62").Value = 2.5

# Add new row 63, copying the date formatting from A62 (so it reuses the
# existing date style instead of creating a new numFmt)
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)

$ws.Range("A63").Value = 45430
$ws.Range("B63").Value = 7
$ws.Range("C63").Formula = "=C62+B63"

$ws.Range("C63").Select() | Out-Null

$wb.Save()
